$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2021-04-17", "overview", "K02000001", "United Kingdom", 4385938, 2206, 35, 127260),
    @("2021-04-18", "overview", "K02000001", "United Kingdom", 4387820, 1882, 10, 127270),
    @("2021-04-19", "overview", "K02000001", "United Kingdom", 4390783, 2963, 4, 127274),
    @("2021-04-20", "overview", "K02000001", "United Kingdom", 4393307, 2524, 33, 127307),
    @("2021-04-21", "overview", "K02000001", "United Kingdom", 4395703, 2396, 22, 127327)
)

$startRow = 249
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    # Column A holds a date-formatted-as-text string (e.g. "2021-04-17").
    # Assigning it straight to .Value would let Excel auto-detect it as a
    # real date and convert it to a serial number, so enter it with a
    # leading apostrophe to force text, then copy the (unstyled) format
    # from an existing plain-text date cell so no extra style is left
    # behind on the new cell.
    $ws.Cells.Item($row, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item(2, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
    $ws.Cells.Item($row, 8).Value = $rowData[7]
}

$excel.CutCopyMode = $false
